$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 854.9048
$ws.Range("I129").Value = 549.125
$ws.Range("J129").Value = 1043.0769
$ws.Range("K129").Value = 1647.375
$ws.Range("L129").Value = 3129.2307
$ws.Range("M129").Value = 3352.625
$ws.Range("N129").Value = -13129.2307

$ws.Range("H137").Value = 4612.212
$ws.Range("I137").Value = 1976.9231
$ws.Range("J137").Value = 6325.15
$ws.Range("K137").Value = 5930.7693
$ws.Range("L137").Value = 18975.45
$ws.Range("M137").Value = -3380.7693
$ws.Range("N137").Value = -24075.45

$ws.Range("H141").Value = 2254.6282
$ws.Range("I141").Value = 1720.4783
$ws.Range("J141").Value = 2478
$ws.Range("K141").Value = 5161.4349
$ws.Range("L141").Value = 7434
$ws.Range("M141").Value = 18.5650999999998
$ws.Range("N141").Value = -17794

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 27000
$ws.Range("J18").Value = 27000
$ws.Range("L18").Value = 27000
$ws.Range("N18").Value = -27644

$ws.Range("H44").Value = 28849
$ws.Range("J44").Value = 28849
$ws.Range("L44").Value = 28849
$ws.Range("N44").Value = -29825

$ws.Range("H55").Value = 27613.666
$ws.Range("I55").Value = 3848
$ws.Range("K55").Value = 3848
$ws.Range("M55").Value = -3533

$ws.Range("H88").Value = 51503.5
$ws.Range("I88").Value = 2906
$ws.Range("J88").Value = 67702.664
$ws.Range("K88").Value = 2906
$ws.Range("L88").Value = 67702.664
$ws.Range("M88").Value = -2500
$ws.Range("N88").Value = -68514.664

$ws.Range("H91").Value = 51503.5
$ws.Range("I91").Value = 2906
$ws.Range("J91").Value = 67702.664
$ws.Range("K91").Value = 2906
$ws.Range("L91").Value = 67702.664
$ws.Range("M91").Value = -1502
$ws.Range("N91").Value = -70510.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 450.5
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -846

$ws.Range("H82").Value = 25937.777
$ws.Range("J82").Value = 35616.91
$ws.Range("L82").Value = 35616.91
$ws.Range("N82").Value = -36382.91

$ws.Range("H85").Value = 25937.777
$ws.Range("J85").Value = 35616.91
$ws.Range("L85").Value = 35616.91
$ws.Range("N85").Value = -38268.91

$ws.Range("H86").Value = 2348.4211
$ws.Range("I86").Value = 2508.1538
$ws.Range("J86").Value = 2002.3334
$ws.Range("K86").Value = 2508.1538
$ws.Range("L86").Value = 2002.3334
$ws.Range("M86").Value = -1385.1538
$ws.Range("N86").Value = -4248.3334

$ws.Range("H89").Value = 2348.4211
$ws.Range("I89").Value = 2508.1538
$ws.Range("J89").Value = 2002.3334
$ws.Range("K89").Value = 12540.769
$ws.Range("L89").Value = 10011.667
$ws.Range("M89").Value = -6924.769
$ws.Range("N89").Value = -21243.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9349.143
$ws.Range("J50").Value = 9349.143
$ws.Range("L50").Value = 9349.143
$ws.Range("N50").Value = -10599.143

$ws.Range("H51").Value = 9230
$ws.Range("J51").Value = 9230
$ws.Range("L51").Value = 9230
$ws.Range("N51").Value = -10702

$ws.Range("H60").Value = 24574.25
$ws.Range("J60").Value = 24574.25
$ws.Range("L60").Value = 24574.25
$ws.Range("N60").Value = -25596.25

$ws.Range("H61").Value = 9230
$ws.Range("J61").Value = 9230
$ws.Range("L61").Value = 9230
$ws.Range("N61").Value = -9926

$ws.Range("H68").Value = 17738.4
$ws.Range("J68").Value = 17738.4
$ws.Range("L68").Value = 17738.4
$ws.Range("N68").Value = -19236.4

$ws.Range("H71").Value = 17738.4
$ws.Range("J71").Value = 17738.4
$ws.Range("L71").Value = 53215.2
$ws.Range("N71").Value = -60703.2

$ws.Range("H109").Value = 14548.571
$ws.Range("J109").Value = 14548.571
$ws.Range("L109").Value = 14548.571
$ws.Range("N109").Value = -16628.571

$ws.Range("H122").Value = 1415.75
$ws.Range("I122").Value = 1332.2858
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3996.8574
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1546.8574
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 66002110
$ws.Range("I81").Value = 600
$ws.Range("J81").Value = 70716504
$ws.Range("K81").Value = 1800
$ws.Range("L81").Value = 212149512
$ws.Range("M81").Value = -677
$ws.Range("N81").Value = -212151758

$ws.Range("H84").Value = 66002110
$ws.Range("I84").Value = 600
$ws.Range("J84").Value = 70716504
$ws.Range("K84").Value = 5400
$ws.Range("L84").Value = 636448536
$ws.Range("M84").Value = 216
$ws.Range("N84").Value = -636459768

$ws.Range("H121").Value = 445610.34
$ws.Range("J121").Value = 703406.8
$ws.Range("L121").Value = 2110220.4
$ws.Range("N121").Value = -2112840.4

$ws.Range("H131").Value = 911.14
$ws.Range("I131").Value = 511.25
$ws.Range("J131").Value = 945.913
$ws.Range("K131").Value = 1533.75
$ws.Range("L131").Value = 2837.739
$ws.Range("M131").Value = 3506.25
$ws.Range("N131").Value = -12917.739

$ws.Range("H140").Value = 21445.346
$ws.Range("I140").Value = 46126.816
$ws.Range("J140").Value = 3345.6
$ws.Range("K140").Value = 138380.448
$ws.Range("L140").Value = 10036.8
$ws.Range("M140").Value = -133200.448
$ws.Range("N140").Value = -20396.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -20970

$ws.Range("H57").Value = 14266.375
$ws.Range("I57").Value = 4870
$ws.Range("K57").Value = 4870
$ws.Range("M57").Value = -4050

$ws.Range("H70").Value = 35645.305
$ws.Range("I70").Value = 39920.207
$ws.Range("J70").Value = 4652.25
$ws.Range("K70").Value = 39920.207
$ws.Range("L70").Value = 4652.25
$ws.Range("M70").Value = -39650.207
$ws.Range("N70").Value = -5192.25

$ws.Range("H73").Value = 35645.305
$ws.Range("I73").Value = 39920.207
$ws.Range("J73").Value = 4652.25
$ws.Range("K73").Value = 39920.207
$ws.Range("L73").Value = 4652.25
$ws.Range("M73").Value = -38984.207
$ws.Range("N73").Value = -6524.25

$ws.Range("H123").Value = 31326
$ws.Range("J123").Value = 31326
$ws.Range("L123").Value = 31326
$ws.Range("N123").Value = -36226

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 923.8421
$ws.Range("I16").Value = 936.2778
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 936.2778
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -766.2778
$ws.Range("N16").Value = -1040

$ws.Range("H122").Value = 5763.636
$ws.Range("I122").Value = 5166.6665
$ws.Range("J122").Value = 6480
$ws.Range("K122").Value = 15499.9995
$ws.Range("L122").Value = 19440
$ws.Range("M122").Value = -13049.9995
$ws.Range("N122").Value = -24340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 19800
$ws.Range("J109").Value = 19800
$ws.Range("L109").Value = 19800
$ws.Range("N109").Value = -22574
